# Refresh the crypto price table: update D/E (price, 1h volume) for the
# existing top-23 coins, then re-write rows 24-51 with the refreshed
# coinranking.com snapshot (a new entry, BitDAO, lands at row 24 and
# pushes the remaining coins down by one row, dropping the last one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Excel COM auto-converts plain-numeric-looking strings (e.g. "15.92")
    # to numbers on assignment. The source sheet stores these as literal
    # text, so force text interpretation, then drop back to the default
    # "Normal" style so no stray number-format is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '30.062.12'
$ws.Range("E2").Value = '  -1.83%  '
# Row 3
$ws.Range("D3").Value = '1.831.50'
$ws.Range("E3").Value = '  -3.09%  '
# Row 4
$ws.Range("E4").Value = '  -0.17%  '
# Row 5
Set-TextValue $ws.Range("D5") '229.11'
$ws.Range("E5").Value = '  -4.13%  '
# Row 6
$ws.Range("E6").Value = '  -0.08%  '
# Row 7
Set-TextValue $ws.Range("D7") '0.4620'
$ws.Range("E7").Value = '  -4.29%  '
# Row 8
Set-TextValue $ws.Range("D8") '0.2687'
$ws.Range("E8").Value = '  -6.12%  '
# Row 9
Set-TextValue $ws.Range("D9") '0.06196'
$ws.Range("E9").Value = '  -5.39%  '
# Row 10
$ws.Range("D10").Value = '1.827.67'
$ws.Range("E10").Value = '  -4.10%  '
# Row 11
Set-TextValue $ws.Range("D11") '0.07339'
$ws.Range("E11").Value = '  -1.77%  '
# Row 12
Set-TextValue $ws.Range("D12") '15.92'
$ws.Range("E12").Value = '  -4.51%  '
# Row 13
Set-TextValue $ws.Range("D13") '4.887'
$ws.Range("E13").Value = '  -4.25%  '
# Row 14
Set-TextValue $ws.Range("D14") '82.50'
$ws.Range("E14").Value = '  -6.38%  '
# Row 15
Set-TextValue $ws.Range("D15") '0.6158'
$ws.Range("E15").Value = '  -7.67%  '
# Row 16
$ws.Range("D16").Value = '30.012.57'
$ws.Range("E16").Value = '  -1.99%  '
# Row 17
$ws.Range("E17").Value = '  -0.12%  '
# Row 18
Set-TextValue $ws.Range("D18") '224.77'
$ws.Range("E18").Value = '  -3.37%  '
# Row 19
Set-TextValue $ws.Range("D19") '0.000007223'
$ws.Range("E19").Value = '  -4.66%  '
# Row 20
Set-TextValue $ws.Range("D20") '0.9975'
$ws.Range("E20").Value = '  -0.48%  '
# Row 21
$ws.Range("D21").Value = '2.070.75'
$ws.Range("E21").Value = '  -2.76%  '
# Row 22
Set-TextValue $ws.Range("D22") '12.23'
$ws.Range("E22").Value = '  -7.92%  '
# Row 23
Set-TextValue $ws.Range("D23") '4.802'
$ws.Range("E23").Value = '  -9.08%  '
# Row 24
$ws.Range("B24").Value = 'BitDAO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
Set-TextValue $ws.Range("D24") '0.3878'
$ws.Range("E24").Value = '  +2.11%  '
# Row 25
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D25") '5.813'
$ws.Range("E25").Value = '  -6.62%  '
# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D26") '164.63'
$ws.Range("E26").Value = '  -2.75%  '
# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D27") '9.079'
$ws.Range("E27").Value = '  -3.02%  '
# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D28") '17.54'
$ws.Range("E28").Value = '  -6.88%  '
# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D29") '1.833'
$ws.Range("E29").Value = '  -6.64%  '
# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D30") '0.1016'
$ws.Range("E30").Value = '  -0.89%  '
# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D31") '1.366'
$ws.Range("E31").Value = '  -2.21%  '
# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D32") '4.039'
$ws.Range("E32").Value = '  -6.70%  '
# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D33") '3.744'
$ws.Range("E33").Value = '  -6.96%  '
# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D34") '0.04767'
$ws.Range("E34").Value = '  -6.05%  '
# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D35") '1.123'
$ws.Range("E35").Value = '  -7.47%  '
# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D36") '0.6924'
$ws.Range("E36").Value = '  -8.17%  '
# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D37") '2.683'
$ws.Range("E37").Value = '  -1.10%  '
# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D38") '0.01809'
$ws.Range("E38").Value = '  -3.57%  '
# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D39") '2.609'
$ws.Range("E39").Value = '  -1.53%  '
# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D40") '0.8818'
$ws.Range("E40").Value = '  -4.24%  '
# Row 41
Set-TextValue $ws.Range("D41") '0.9999'
$ws.Range("E41").Value = '  -0.26%  '
# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D42") '1.901'
$ws.Range("E42").Value = '  -8.20%  '
# Row 43
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D43") '102.90'
$ws.Range("E43").Value = '  -4.13%  '
# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D44") '5.432'
$ws.Range("E44").Value = '  -4.07%  '
# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D45") '0.3967'
$ws.Range("E45").Value = '  -7.72%  '
# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D46") '6.850'
$ws.Range("E46").Value = '  -7.77%  '
# Row 47
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D47") '0.1179'
$ws.Range("E47").Value = '  -7.38%  '
# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D48") '58.82'
$ws.Range("E48").Value = '  -8.83%  '
# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D49") '8.425'
$ws.Range("E49").Value = '  -6.31%  '
# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.05525'
$ws.Range("E50").Value = '  -2.53%  '
# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range("D51") '32.41'
$ws.Range("E51").Value = '  -4.44%  '
